# Auto-generated: refresh market-price derived columns (H-N) across all sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 11494428
$ws.Cells.Item(6, 9).Value = 11494428
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 34483284
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -34483172
$ws.Cells.Item(6, 14).Value = ""
$ws.Cells.Item(9, 8).Value = 356.2
$ws.Cells.Item(9, 10).Value = 1611
$ws.Cells.Item(9, 12).Value = 1611
$ws.Cells.Item(9, 14).Value = -1949
$ws.Cells.Item(28, 8).Value = 4650.923
$ws.Cells.Item(28, 9).Value = 2838
$ws.Cells.Item(28, 10).Value = 8730
$ws.Cells.Item(28, 11).Value = 2838
$ws.Cells.Item(28, 12).Value = 8730
$ws.Cells.Item(28, 13).Value = -2353
$ws.Cells.Item(28, 14).Value = -9700
$ws.Cells.Item(129, 8).Value = 1552.5454
$ws.Cells.Item(129, 9).Value = 429.8
$ws.Cells.Item(129, 10).Value = 2488.1667
$ws.Cells.Item(129, 11).Value = 1289.4
$ws.Cells.Item(129, 12).Value = 7464.500100000001
$ws.Cells.Item(129, 13).Value = 3710.6
$ws.Cells.Item(129, 14).Value = -17464.5001
$ws.Cells.Item(137, 8).Value = 2334.5
$ws.Cells.Item(137, 9).Value = 1849.375
$ws.Cells.Item(137, 10).Value = 2765.7222
$ws.Cells.Item(137, 11).Value = 5548.125
$ws.Cells.Item(137, 12).Value = 8297.1666
$ws.Cells.Item(137, 13).Value = -2998.125
$ws.Cells.Item(137, 14).Value = -13397.1666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3355.2778
$ws.Cells.Item(2, 9).Value = 1773.909
$ws.Cells.Item(2, 10).Value = 5840.2856
$ws.Cells.Item(2, 11).Value = 1773.909
$ws.Cells.Item(2, 12).Value = 5840.2856
$ws.Cells.Item(2, 13).Value = -1660.909
$ws.Cells.Item(2, 14).Value = -6066.2856
$ws.Cells.Item(32, 8).Value = 1509133.5
$ws.Cells.Item(32, 9).Value = 1788447.9
$ws.Cells.Item(32, 10).Value = 5132.615
$ws.Cells.Item(32, 11).Value = 1788447.9
$ws.Cells.Item(32, 12).Value = 5132.615
$ws.Cells.Item(32, 13).Value = -1788160.9
$ws.Cells.Item(32, 14).Value = -5706.615
$ws.Cells.Item(61, 8).Value = 7384.6924
$ws.Cells.Item(61, 9).Value = 3391.2273
$ws.Cells.Item(61, 10).Value = 12552.706
$ws.Cells.Item(61, 11).Value = 3391.2273
$ws.Cells.Item(61, 12).Value = 12552.706
$ws.Cells.Item(61, 13).Value = -3179.2273
$ws.Cells.Item(61, 14).Value = -12976.706
$ws.Cells.Item(116, 8).Value = 3355.2778
$ws.Cells.Item(116, 9).Value = 1773.909
$ws.Cells.Item(116, 10).Value = 5840.2856
$ws.Cells.Item(116, 11).Value = 1773.909
$ws.Cells.Item(116, 12).Value = 5840.2856
$ws.Cells.Item(116, 13).Value = 520.0909999999999
$ws.Cells.Item(116, 14).Value = -10428.2856
$ws.Cells.Item(136, 8).Value = 7384.6924
$ws.Cells.Item(136, 9).Value = 3391.2273
$ws.Cells.Item(136, 10).Value = 12552.706
$ws.Cells.Item(136, 11).Value = 10173.6819
$ws.Cells.Item(136, 12).Value = 37658.118
$ws.Cells.Item(136, 13).Value = -7623.6819
$ws.Cells.Item(136, 14).Value = -42758.118

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3355.2778
$ws.Cells.Item(3, 9).Value = 1773.909
$ws.Cells.Item(3, 10).Value = 5840.2856
$ws.Cells.Item(3, 11).Value = 1773.909
$ws.Cells.Item(3, 12).Value = 5840.2856
$ws.Cells.Item(3, 13).Value = -1659.909
$ws.Cells.Item(3, 14).Value = -6068.2856
$ws.Cells.Item(20, 8).Value = 12822503
$ws.Cells.Item(20, 9).Value = 20834736
$ws.Cells.Item(20, 11).Value = 20834736
$ws.Cells.Item(20, 13).Value = -20834489
$ws.Cells.Item(105, 8).Value = 3675.96
$ws.Cells.Item(105, 9).Value = 2507.4167
$ws.Cells.Item(105, 11).Value = 2507.4167
$ws.Cells.Item(105, 13).Value = -760.4167000000002
$ws.Cells.Item(134, 8).Value = 6964.9
$ws.Cells.Item(134, 9).Value = 1015.2222
$ws.Cells.Item(134, 11).Value = 3045.6666
$ws.Cells.Item(134, 13).Value = -510.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6060.8076
$ws.Cells.Item(31, 9).Value = 2920.4412
$ws.Cells.Item(31, 10).Value = 11992.611
$ws.Cells.Item(31, 11).Value = 2920.4412
$ws.Cells.Item(31, 12).Value = 11992.611
$ws.Cells.Item(31, 13).Value = -2625.4412
$ws.Cells.Item(31, 14).Value = -12582.611
$ws.Cells.Item(34, 8).Value = 6060.8076
$ws.Cells.Item(34, 9).Value = 2920.4412
$ws.Cells.Item(34, 10).Value = 11992.611
$ws.Cells.Item(34, 11).Value = 2920.4412
$ws.Cells.Item(34, 12).Value = 11992.611
$ws.Cells.Item(34, 13).Value = -2718.4412
$ws.Cells.Item(34, 14).Value = -12396.611
$ws.Cells.Item(132, 8).Value = 4098.5625
$ws.Cells.Item(132, 9).Value = 2215.8
$ws.Cells.Item(132, 10).Value = 7236.5
$ws.Cells.Item(132, 11).Value = 6647.400000000001
$ws.Cells.Item(132, 12).Value = 21709.5
$ws.Cells.Item(132, 13).Value = -4117.400000000001
$ws.Cells.Item(132, 14).Value = -26769.5
$ws.Cells.Item(134, 8).Value = 3765.8596
$ws.Cells.Item(134, 9).Value = 2125.9565
$ws.Cells.Item(134, 11).Value = 6377.869499999999
$ws.Cells.Item(134, 13).Value = -3842.869499999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1113620
$ws.Cells.Item(5, 10).Value = 2896
$ws.Cells.Item(5, 12).Value = 8688
$ws.Cells.Item(5, 14).Value = -8912
$ws.Cells.Item(41, 8).Value = 1695
$ws.Cells.Item(41, 10).Value = 1875
$ws.Cells.Item(41, 12).Value = 5625
$ws.Cells.Item(41, 14).Value = -6301
$ws.Cells.Item(68, 8).Value = 2347.7693
$ws.Cells.Item(68, 10).Value = 2350.4243
$ws.Cells.Item(68, 12).Value = 7051.2729
$ws.Cells.Item(68, 14).Value = -8673.2729
$ws.Cells.Item(71, 8).Value = 2347.7693
$ws.Cells.Item(71, 10).Value = 2350.4243
$ws.Cells.Item(71, 12).Value = 21153.8187
$ws.Cells.Item(71, 14).Value = -29265.8187
$ws.Cells.Item(87, 8).Value = 11970.588
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 11970.588
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 35911.764
$ws.Cells.Item(87, 13).Value = ""
$ws.Cells.Item(87, 14).Value = -38407.764
$ws.Cells.Item(90, 8).Value = 11970.588
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 11970.588
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 12).Value = 107735.292
$ws.Cells.Item(90, 13).Value = ""
$ws.Cells.Item(90, 14).Value = -120215.292
$ws.Cells.Item(122, 8).Value = 1416041.2
$ws.Cells.Item(122, 9).Value = 2830088
$ws.Cells.Item(122, 10).Value = 1994.5
$ws.Cells.Item(122, 11).Value = 25470792
$ws.Cells.Item(122, 12).Value = 17950.5
$ws.Cells.Item(122, 13).Value = -25468342
$ws.Cells.Item(122, 14).Value = -22850.5
$ws.Cells.Item(135, 8).Value = 1113620
$ws.Cells.Item(135, 10).Value = 2896
$ws.Cells.Item(135, 12).Value = 26064
$ws.Cells.Item(135, 14).Value = -31134
$ws.Cells.Item(140, 8).Value = 119048.35
$ws.Cells.Item(140, 9).Value = 138922.23
$ws.Cells.Item(140, 11).Value = 416766.6900000001
$ws.Cells.Item(140, 13).Value = -411586.6900000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 3333511.8
$ws.Cells.Item(2, 9).Value = 145.25
$ws.Cells.Item(2, 11).Value = 145.25
$ws.Cells.Item(2, 13).Value = -32.25
$ws.Cells.Item(52, 8).Value = 86258
$ws.Cells.Item(52, 10).Value = 86258
$ws.Cells.Item(52, 12).Value = 86258
$ws.Cells.Item(52, 14).Value = -86776
$ws.Cells.Item(97, 8).Value = 2755.8
$ws.Cells.Item(97, 9).Value = 2968.375
$ws.Cells.Item(97, 11).Value = 2968.375
$ws.Cells.Item(97, 13).Value = -2472.375
$ws.Cells.Item(102, 8).Value = 2928.9167
$ws.Cells.Item(102, 9).Value = 1555.875
$ws.Cells.Item(102, 11).Value = 1555.875
$ws.Cells.Item(102, 13).Value = 66.125
$ws.Cells.Item(113, 8).Value = 7141.9287
$ws.Cells.Item(113, 9).Value = 4944.4287
$ws.Cells.Item(113, 10).Value = 7874.4287
$ws.Cells.Item(113, 11).Value = 4944.4287
$ws.Cells.Item(113, 12).Value = 7874.4287
$ws.Cells.Item(113, 13).Value = -2774.4287
$ws.Cells.Item(113, 14).Value = -12214.4287
$ws.Cells.Item(122, 8).Value = 25670000
$ws.Cells.Item(122, 9).Value = 33369224
$ws.Cells.Item(122, 11).Value = 100107672
$ws.Cells.Item(122, 13).Value = -100105222
$ws.Cells.Item(132, 8).Value = 6402.7856
$ws.Cells.Item(132, 10).Value = 8984
$ws.Cells.Item(132, 12).Value = 26952
$ws.Cells.Item(132, 14).Value = -32012

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7148.36
$ws.Cells.Item(7, 9).Value = 4522.1113
$ws.Cells.Item(7, 10).Value = 8625.625
$ws.Cells.Item(7, 11).Value = 4522.1113
$ws.Cells.Item(7, 12).Value = 8625.625
$ws.Cells.Item(7, 13).Value = -4410.1113
$ws.Cells.Item(7, 14).Value = -8849.625
$ws.Cells.Item(16, 8).Value = 3483
$ws.Cells.Item(16, 9).Value = 3483
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 3483
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -3313
$ws.Cells.Item(16, 14).Value = ""
$ws.Cells.Item(122, 8).Value = 7158.643
$ws.Cells.Item(122, 9).Value = 7469.148
$ws.Cells.Item(122, 10).Value = 6599.7334
$ws.Cells.Item(122, 11).Value = 22407.444
$ws.Cells.Item(122, 12).Value = 19799.2002
$ws.Cells.Item(122, 13).Value = -19957.444
$ws.Cells.Item(122, 14).Value = -24699.2002
$ws.Cells.Item(126, 8).Value = 7148.36
$ws.Cells.Item(126, 9).Value = 4522.1113
$ws.Cells.Item(126, 10).Value = 8625.625
$ws.Cells.Item(126, 11).Value = 13566.3339
$ws.Cells.Item(126, 12).Value = 25876.875
$ws.Cells.Item(126, 13).Value = -11096.3339
$ws.Cells.Item(126, 14).Value = -30816.875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 12212188
$ws.Cells.Item(132, 9).Value = 21754740
$ws.Cells.Item(132, 10).Value = 18927.611
$ws.Cells.Item(132, 11).Value = 65264220
$ws.Cells.Item(132, 12).Value = 56782.833
$ws.Cells.Item(132, 13).Value = -65261690
$ws.Cells.Item(132, 14).Value = -61842.833
$ws.Cells.Item(136, 8).Value = 58888664
$ws.Cells.Item(136, 9).Value = 142859310
$ws.Cells.Item(136, 10).Value = 109211
$ws.Cells.Item(136, 11).Value = 428577930
$ws.Cells.Item(136, 12).Value = 327633
$ws.Cells.Item(136, 13).Value = -428575380
$ws.Cells.Item(136, 14).Value = -332733
